# "deleted pics and added readme" -- apply the new Trials/Stimuli list to
# Tabelle1 (the file actually adds a new stimulus row: Stimuli/9230.jpg,
# inserted in numeric-sorted order, pushing the following rows down).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Final, full list of image files that should occupy A2:A21 after the edit
# (A1 keeps the "ImageFile" header).
$values = @(
    "Stimuli/227.jpg",
    "Stimuli/252.jpg",
    "Stimuli/1051.jpg",
    "Stimuli/2800.jpg",
    "Stimuli/3061.jpg",
    "Stimuli/3230.jpg",
    "Stimuli/6561.jpg",
    "Stimuli/6838.jpg",
    "Stimuli/9120.jpg",
    "Stimuli/9181.jpg",
    "Stimuli/9185.jpg",
    "Stimuli/9230.jpg",
    "Stimuli/9254.jpg",
    "Stimuli/9295.jpg",
    "Stimuli/9332.jpg",
    "Stimuli/9411.jpg",
    "Stimuli/9420.jpg",
    "Stimuli/9421.jpg",
    "Stimuli/9599.jpg",
    "Stimuli/9905.jpg"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("G15").Select()
